# Applies the "SegmentedSieve don't have 2,3 in Prime no" edit to the
# DSA Supreme List worksheet.
#
# Summary of the change:
#  - Row 71 (Count Primes): label gets a clarifying suffix
#       "Count Primes (204)" -> "Count Primes (204)   [using Sieve of Eratosthenes]"
#    (the T.C. and LINK values for that row stay the same)
#  - Row 75 (previously just had the index "6" in col A) gets filled in with a
#    new entry: "Prime No in Range using Segmented Sieve"
#  - A brand new row 76 is appended with index "7": "Product of Primes"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 71: clarify the Count Primes question title -----------------------
$ws.Range("B71").Value = "Count Primes (204)   [using Sieve of Eratosthenes]"

# --- Row 75: new "Segmented Sieve" entry ------------------------------------
# Copy formatting from the row above (row 74) down onto row 75 first, so the
# new row matches the rest of the table's look (fonts/alignment/number
# formats), then fill in the values.
$ws.Range("A74:D74").Copy()
$ws.Range("A75:D75").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A75").Value = 6
$ws.Range("B75").Value = "Prime No in Range using Segmented Sieve"
$ws.Range("C75").Value = "n(log(logn))"
$ws.Range("D75").Value = "https://github.com/rohillanishant/DSA-Cpp/blob/master/SieveOfEratosthenes/SegmentedSieve.cpp"
$ws.Rows("75").RowHeight = 15.6

# --- Row 76: brand new "Product of Primes" entry ----------------------------
$ws.Range("A75:D75").Copy()
$ws.Range("A76:D76").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Product of Primes"
$ws.Range("C76").Value = "n(log(logn))"
$ws.Range("D76").Value = "https://practice.geeksforgeeks.org/problems/product-of-primes5328/1"
$ws.Rows("76").RowHeight = 15.6

# --- Fix up selection / view state like the real edit session --------------
$ws.Range("D63").Select()
